$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update year-period header labels (shift one year forward) ---
$ws.Range("E8").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E27").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F27").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G27").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H27").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I27").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E35").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F35").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G35").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H35").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I35").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E43").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F43").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G43").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H43").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I43").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E51").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F51").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G51").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H51").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I51").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E59").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F59").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G59").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H59").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I59").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E67").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F67").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G67").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H67").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I67").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E75").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F75").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G75").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H75").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I75").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E83").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F83").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G83").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H83").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I83").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E91").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F91").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G91").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H91").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I91").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E98").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F98").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G98").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H98").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I98").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E105").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F105").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G105").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H105").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I105").Value = "دوازده ماهه منتهی به 1401/12"

# --- Shift yearly data columns left (drop oldest year, add new year) ---
$ws.Range("E10").Value = 6139276
$ws.Range("F10").Value = 6782515
$ws.Range("G10").Value = 13672872
$ws.Range("H10").Value = 45731705
$ws.Range("I10").Value = 46464713
$ws.Range("E11").Value = 380948
$ws.Range("F11").Value = 459138
$ws.Range("G11").Value = 1617607
$ws.Range("H11").Value = 2234393
$ws.Range("I11").Value = 3816668
$ws.Range("E12").Value = 3717771
$ws.Range("F12").Value = 4489574
$ws.Range("G12").Value = 5724950
$ws.Range("H12").Value = 12277785
$ws.Range("I12").Value = 16454698
$ws.Range("E13").Value = 10237995
$ws.Range("F13").Value = 11731227
$ws.Range("G13").Value = 21015429
$ws.Range("H13").Value = 60243883
$ws.Range("I13").Value = 66736079
$ws.Range("E15").Value = 10237995
$ws.Range("F15").Value = 11731227
$ws.Range("G15").Value = 21015429
$ws.Range("H15").Value = 60243883
$ws.Range("I15").Value = 66736079
$ws.Range("E18").Value = 10237995
$ws.Range("F18").Value = 11731227
$ws.Range("G18").Value = 21015429
$ws.Range("H18").Value = 60243883
$ws.Range("I18").Value = 66736079
$ws.Range("E19").Value = 832101
$ws.Range("F19").Value = 1347732
$ws.Range("G19").Value = 1787555
$ws.Range("H19").Value = 4233413
$ws.Range("I19").Value = 9404229
$ws.Range("E20").Value = -1347732
$ws.Range("F20").Value = -1787555
$ws.Range("G20").Value = -4233413
$ws.Range("H20").Value = -9404229
$ws.Range("I20").Value = -9181565
$ws.Range("E21").Value = 9722364
$ws.Range("F21").Value = 11291404
$ws.Range("G21").Value = 18569571
$ws.Range("H21").Value = 55073067
$ws.Range("I21").Value = 66958743
$ws.Range("E23").Value = 9722364
$ws.Range("F23").Value = 11291404
$ws.Range("G23").Value = 18569571
$ws.Range("H23").Value = 55073067
$ws.Range("I23").Value = 66958743
$ws.Range("H29").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("H30").Value = "-"
$ws.Range("H37").Value = 763
$ws.Range("I37").Value = 774
$ws.Range("E38").Value = 25904
$ws.Range("F38").Value = 24253
$ws.Range("G38").Value = 26004
$ws.Range("H38").Value = "-"
$ws.Range("E39").Value = 25904
$ws.Range("F39").Value = 24253
$ws.Range("G39").Value = 26004
$ws.Range("H39").Value = 763
$ws.Range("I39").Value = 774
$ws.Range("H45").Value = 763
$ws.Range("I45").Value = 774
$ws.Range("E46").Value = 25904
$ws.Range("F46").Value = 24253
$ws.Range("G46").Value = 26004
$ws.Range("H46").Value = "-"
$ws.Range("E47").Value = 25904
$ws.Range("F47").Value = 24253
$ws.Range("G47").Value = 26004
$ws.Range("H47").Value = 763
$ws.Range("I47").Value = 774
$ws.Range("H53").Value = 0
$ws.Range("E54").Value = 0
$ws.Range("H54").Value = "-"
$ws.Range("H61").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("H62").Value = "-"
$ws.Range("H69").Value = 45731705
$ws.Range("I69").Value = 46464713
$ws.Range("E70").Value = 6139276
$ws.Range("F70").Value = 6782515
$ws.Range("G70").Value = 13672872
$ws.Range("H70").Value = "-"
$ws.Range("E71").Value = 6139276
$ws.Range("F71").Value = 6782515
$ws.Range("G71").Value = 13672872
$ws.Range("H71").Value = 45731705
$ws.Range("I71").Value = 46464713
$ws.Range("H77").Value = 45731705
$ws.Range("I77").Value = 46464713
$ws.Range("E78").Value = 6139276
$ws.Range("F78").Value = 6782515
$ws.Range("G78").Value = 13672872
$ws.Range("H78").Value = "-"
$ws.Range("E79").Value = 6139276
$ws.Range("F79").Value = 6782515
$ws.Range("G79").Value = 13672872
$ws.Range("H79").Value = 45731705
$ws.Range("I79").Value = 46464713
$ws.Range("H85").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("H86").Value = "-"
$ws.Range("H93").Value = 59936703801
$ws.Range("I93").Value = 60031928941
$ws.Range("E94").Value = 237000715
$ws.Range("F94").Value = 279656744
$ws.Range("G94").Value = 525798800
$ws.Range("H94").Value = "-"
$ws.Range("H100").Value = 59936703801
$ws.Range("I100").Value = 60031928941
$ws.Range("E101").Value = 237000715
$ws.Range("F101").Value = 279656744
$ws.Range("G101").Value = 525798800
$ws.Range("H101").Value = "-"
$ws.Range("E107").Value = 86744
$ws.Range("F107").Value = 107058
$ws.Range("G107").Value = 0
$ws.Range("E111").Value = 638191
$ws.Range("F111").Value = 812293
$ws.Range("G111").Value = 1584414
$ws.Range("H111").Value = 1738726
$ws.Range("I111").Value = 1754775
$ws.Range("E112").Value = 675637
$ws.Range("F112").Value = 1057132
$ws.Range("G112").Value = 990806
$ws.Range("H112").Value = 5646581
$ws.Range("I112").Value = 7909282
$ws.Range("E113").Value = 177455
$ws.Range("F113").Value = 135249
$ws.Range("G113").Value = 400424
$ws.Range("H113").Value = 484188
$ws.Range("I113").Value = 651994
$ws.Range("E114").Value = 1621314
$ws.Range("F114").Value = 1569205
$ws.Range("G114").Value = 2054464
$ws.Range("H114").Value = 2851858
$ws.Range("I114").Value = 4777538
$ws.Range("E116").Value = 518430
$ws.Range("F116").Value = 808637
$ws.Range("G116").Value = 694842
$ws.Range("H116").Value = 1556432
$ws.Range("I116").Value = 1361109
$ws.Range("E117").Value = 3717771
$ws.Range("F117").Value = 4489574
$ws.Range("G117").Value = 5724950
$ws.Range("H117").Value = 12277785
$ws.Range("I117").Value = 16454698
